# Apply the "5.0.0 -> 6.0.0 / v3-PatientImportance" metadata refresh to the
# ValueSet-patient-importance workbook.

$wb = $excel.ActiveWorkbook

# --- 1. Rename the second sheet -------------------------------------------
$wsInclude = $wb.Worksheets.Item(2)
$wsInclude.Name = "Include ValueSets"

# --- 2. Metadata sheet (sheet 1) updates -----------------------------------
$wsMeta = $wb.Worksheets.Item(1)

# Version: 5.0.0 -> 6.0.0
$wsMeta.Range("B3").Value = "6.0.0"

# Date: refresh the publish timestamp
$wsMeta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value (was blank) -> "Alvearie Team"
$wsMeta.Range("B9").Value = "Alvearie Team"

# Row 10 used to be a second "Contact" row -> becomes "Jurisdiction"
$wsMeta.Range("A10").Value = "Jurisdiction"
$wsMeta.Range("B10").Value = "United States of America"

# Row 11 used to be a third "Contact" row -> becomes "Description"
# (takes the description text that used to live on row 12)
$wsMeta.Range("A11").Value = "Description"
$wsMeta.Range("B11").Value = "Patient importance status codes derived from customer-specific code mappings, used to trigger or limit IBM Health Data Connect patient operations."

# Row 12 used to be "Description" -> becomes "Purpose" (value cleared)
$wsMeta.Range("A12").Value = "Purpose"
$wsMeta.Range("B12").ClearContents()

# Row 13 used to be "Purpose" -> becomes "Copyright" (value stays empty)
$wsMeta.Range("A13").Value = "Copyright"

# Row 14 used to be "Copyright" -> becomes "Immutable" with its value
$wsMeta.Range("A14").Value = "Immutable"
$wsMeta.Range("B14").Value = "BooleanType[null]"

# Row 15 (old "Immutable"/"BooleanType[null]") is no longer needed now that
# row 14 carries it - delete the row entirely so the sheet shrinks to B14.
$wsMeta.Rows.Item(15).Delete()

# --- 3. "Include ValueSets" sheet (sheet 2) updates -------------------------
$wsInclude.Range("A1").Value = "ValueSet URL"
$wsInclude.Range("A2").Value = "http://terminology.hl7.org/ValueSet/v3-PatientImportance"

# Old rows 3-4 (and their column-B data) are dropped entirely.
$wsInclude.Range("A3:A4").EntireRow.Delete()
